# Update the "Förändrad" date column (C2:C14) from 45243 (2023-11-13)
# to 45244 (2023-11-14), as produced by the automatic update process.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
